# Apply updated cryptocurrency market data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.003.44'
$ws.Range('E2').Value = '  +5.47%  '
$ws.Range('D3').Value = '2.362.41'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('D5').Value = "'548.84"
$ws.Range('E5').Value = '  +2.92%  '
$ws.Range('D6').Value = "'132.96"
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').Value = '2.360.00'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('D11').Value = "'5.51"
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = "'24.09"
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('D15').Value = '2.782.65'
$ws.Range('E15').Value = '  +3.14%  '
$ws.Range('D16').Value = '60.898.47'
$ws.Range('E16').Value = '  +5.37%  '
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').Value = '2.375.95'
$ws.Range('E18').Value = '  +3.67%  '
$ws.Range('D19').Value = "'10.75"
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = "'6.89"
$ws.Range('E21').Value = '  +8.39%  '
$ws.Range('D22').Value = "'315.93"
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = "'63.54"
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('E25').Value = '  +4.69%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = "'8.02"
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +6.64%  '
$ws.Range('D29').Value = "'1.77"
$ws.Range('E29').Value = '  +3.78%  '
$ws.Range('D30').Value = "'172.31"
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('D31').Value = '0.0₃0737'
$ws.Range('E31').Value = '  +3.03%  '
$ws.Range('D32').Value = "'1.15"
$ws.Range('E32').Value = '  +10.24%  '
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').Value = "'1.43"
$ws.Range('E34').Value = '  +16.10%  '
$ws.Range('D35').Value = "'0.384"
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D36').Value = "'18.09"
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').Value = "'4.18"
$ws.Range('E39').Value = '  +7.61%  '
$ws.Range('D40').Value = "'317.00"
$ws.Range('E40').Value = '  +9.73%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'1.54"
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = "'38.30"
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').Value = "'143.56"
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('D45').Value = "'0.0957"
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = "'19.49"
$ws.Range('E46').Value = '  +7.71%  '
$ws.Range('D47').Value = "'0.0501"
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').Value = "'0.565"
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('D50').Value = '0.0₆0212'
$ws.Range('E50').Value = '  +6.02%  '
$ws.Range('E51').Value = '  +0.95%  '
